# Slide 3 ("Presentation, Code Samples and Resources") - Content Placeholder 4
# Update the "Custom Angular Module" bullet + insert new logging bullets +
# merge the "Guide (markdown/PDF)" bullet back into a single run.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(3)
$sh = $s.Shapes.Item(2)   # "Content Placeholder 4"

# --- 1) "Custom Angular Module" -> "Custom Angular Module (source, dist)"
#        split across 3 runs to match authored formatting.
$tr   = $sh.TextFrame.TextRange
$para = $tr.Paragraphs(2, 1)
$run1 = $para.Runs(1, 1)
$run1.Text = "Custom Angular Module (source, "

$tr   = $sh.TextFrame.TextRange
$para = $tr.Paragraphs(2, 1)
$para.InsertAfter("dist")

$tr   = $sh.TextFrame.TextRange
$para = $tr.Paragraphs(2, 1)
$para.InsertAfter(")")

# --- 2) After "(2) client Angular applications ..." insert two new
#        sub-bullets: "Basic logging" and "Logging with configuration".
$tr   = $sh.TextFrame.TextRange
$para = $tr.Paragraphs(3, 1)
$para.InsertAfter("`rBasic logging`rLogging ")

$tr   = $sh.TextFrame.TextRange
$para = $tr.Paragraphs(5, 1)
$para.InsertAfter("with configuration")

$tr   = $sh.TextFrame.TextRange
$para = $tr.Paragraphs(4, 1)
$para.IndentLevel = 3

$tr   = $sh.TextFrame.TextRange
$para = $tr.Paragraphs(5, 1)
$para.IndentLevel = 3

# --- 3) "Guide " + "(markdown/PDF)" -> single run "Guide (markdown/PDF)"
$tr   = $sh.TextFrame.TextRange
$para = $tr.Paragraphs(7, 1)
$run1 = $para.Runs(1, 1)
$run1.Text = "Guide (markdown/PDF)"

$tr   = $sh.TextFrame.TextRange
$para = $tr.Paragraphs(7, 1)
$run2 = $para.Runs(2, 1)
$run2.Text = ""
